# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# worksheets to match the newly generated gh-pages data output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 8597
$ws1.Range("F3").Value  = 73
$ws1.Range("F4").Value  = 30
$ws1.Range("F5").Value  = 86
$ws1.Range("F6").Value  = 1364
$ws1.Range("F7").Value  = 122
$ws1.Range("F9").Value  = 31
$ws1.Range("F10").Value = 9312
$ws1.Range("F12").Value = 95
$ws1.Range("F13").Value = 218
$ws1.Range("F14").Value = 176
$ws1.Range("F15").Value = 353
$ws1.Range("F16").Value = 6272
$ws1.Range("F17").Value = 1058
$ws1.Range("F18").Value = 76
$ws1.Range("F19").Value = 42
$ws1.Range("F20").Value = 124

# --- Sheet "全部类型" --------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 8597
$ws4.Range("F3").Value  = 73
$ws4.Range("F4").Value  = 30
$ws4.Range("F5").Value  = 86
$ws4.Range("F6").Value  = 1364
$ws4.Range("F7").Value  = 122
$ws4.Range("F9").Value  = 31
$ws4.Range("F12").Value = 9312
$ws4.Range("F14").Value = 95
$ws4.Range("F15").Value = 218
$ws4.Range("F16").Value = 176
$ws4.Range("F17").Value = 353
$ws4.Range("F18").Value = 6272
$ws4.Range("F19").Value = 1058
$ws4.Range("F20").Value = 76
$ws4.Range("F21").Value = 43
$ws4.Range("F22").Value = 124

$wb.Save()
